# #12 Bubbles colors updated on solution slides on presentation
#
# 1) Bump the deck's cached "datetimeFigureOut" placeholder text on the
#    slide master and every slide layout (4/5/2020 -> 4/6/2020).
# 2) Recolor the three speech-bubble callout shapes on slide 18 (the
#    "Result" slide) with a light-blue fill/outline (00B0F0).

$p = $ppt.ActivePresentation

# --- 1. Update the cached date field text wherever it appears --------------
$oldDate = "4/5/2020"
$newDate = "4/6/2020"

function Update-DatePlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholderText $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholderText $layout.Shapes
}

# --- 2. Recolor the narrative bubble callouts on slide 18 ------------------
# PowerPoint's OLE color values are packed as 0x00BBGGRR, so build it from
# the target hex color 00B0F0 (R=0x00, G=0xB0, B=0xF0) ourselves since this
# host does not expose the VBA RGB() helper.
function Get-OleColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}
$bubbleColor = Get-OleColor 0x00 0xB0 0xF0  # &H00B0F0 -> 15773696

$slide = $p.Slides.Item(18)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -like "Bulle narrative*") {
        $shp.Fill.ForeColor.RGB = $bubbleColor
        $shp.Line.ForeColor.RGB = $bubbleColor
    }
}
